$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out row 4 (column A no longer used; B4/C4 values cleared but keep style)
$ws.Cells.Item(4, 1).ClearContents()
$ws.Cells.Item(4, 2).ClearContents()
$ws.Cells.Item(4, 3).ClearContents()

# Extend formatting (style "2" = time number format) down through row 25 for columns B and C
$ws.Range("B4:C4").Copy() | Out-Null
$ws.Range("B5:C25").PasteSpecial(-4122) | Out-Null

# Update the active selection to match target state
$ws.Range("C16").Select() | Out-Null
